$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: the call-off order becomes confirmed, and gets assigned the
# required delivery date / quantity values (copied from the call-off's
# current delivery date and ordered quantity).
$ws.Range("M2").Value2 = $ws.Range("M2").Value2 + 13
$ws.Range("O2").Value2 = "confirmed"
$ws.Range("S2").NumberFormat = $ws.Range("M2").NumberFormat
$ws.Range("S2").Value2 = $ws.Range("M2").Value2
$ws.Range("T2").Value2 = $ws.Range("K2").Value2

# Row 3: remaining (un-confirmed) part of the order keeps its own
# postfix bumped, and the quantity that was moved to the confirmed
# part above is subtracted, leaving the remainder for this row.
$ws.Range("D3").Value2 = "2437589/9"
$ws.Range("K3").Value2 = $ws.Range("K3").Value2 - $ws.Range("K2").Value2
$ws.Range("L3").Value2 = $ws.Range("L3").Value2 - $ws.Range("L2").Value2
